# This script applies a data correction to the "Mexico Liga MX Femenil" sheet.
# For a number of pairs of adjacent rows, the match-record data in columns
# B (id) through AC (closing odds) was entered in swapped order. Column A
# (the sequential row index) stays attached to its row. We fix this by
# swapping the B:AC values between each pair of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Mexico Liga MX Femenil")

# Row pairs (1-based worksheet row numbers) whose B:AC contents must be swapped.
$pairs = @(
    @(71, 72),
    @(101, 102),
    @(131, 132),
    @(133, 134),
    @(213, 214),
    @(215, 216),
    @(229, 230),
    @(232, 233),
    @(248, 249),
    @(251, 252),
    @(271, 272)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B${r1}:AC${r1}")
    $range2 = $ws.Range("B${r2}:AC${r2}")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
